$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = -3977691690.71
$ws.Range("P2").Value = 2243545641.88
$ws.Range("Q2").Value = 2862183889.56
$ws.Range("R2").Value = -23.0415046616
$ws.Range("S2").Value = 1607583269.32
$ws.Range("T2").Value = 1607583269.32
$ws.Range("U2").Value = -25.1326749054
$ws.Range("V2").Value = 379786578.71
$ws.Range("W2").Value = 474763697.82
$ws.Range("X2").Value = 323002627.82
$ws.Range("Y2").Value = -3814499435.1
$ws.Range("Z2").Value = -3764160557.47
$ws.Range("AA2").Value = 276285937.88
$ws.Range("AG2").Value = 13074207.97
$ws.Range("AP2").Value = -41.183271055
$ws.Range("AQ2").Value = -307.800190863374
$ws.Range("AR2").Value = -240.363970872786
$ws.Range("AS2").Value = -3982881065.78
$ws.Range("AT2").Value = -235.697907242106
